$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$rows = @(
    @("2026-02-01", "15:56:42", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "15:56:42", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "15:56:49", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "15:57:00", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "15:57:10", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "15:57:20", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "15:57:30", "15:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "15:57:41", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 30
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    # Column A holds a date-like string ("2026-02-01") which Excel would
    # otherwise auto-convert into a date serial number. Force it to be
    # stored as text (matching the rest of the log), then clear the
    # number-format override so no extra style gets attached to the cell.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row[0]
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
